{"js": "// The paragraph currently contains a Word field \"{ m:enduserdoc }\" (a\n// fldChar begin / instrText / fldChar end run triple) followed by a bold\n// red run with an \"unexpected tag\" error message. The M2Doc parser was\n// switched to a token-iterator based field rewriter, so the expected\n// generation output now shows the raw field code as literal template\n// text (\"{m:enduserdoc}\") instead of an actual Word field, and the error\n// message gets a leading arrow marker (\"    <---\").\n\nconst body = context.document.body;\n\n// Step 1: insert the literal text \"{m:enduserdoc}\" at the very start of\n// the (only) paragraph, i.e. right before the field. Doing this before\n// removing the field means the new text lands next to the field's begin\n// run, so it does not pick up the bold / red formatting of the trailing\n// error-message run.\nconst firstParagraph = body.paragraphs.getFirst();\nconst startRange = firstParagraph.getRange(\"Start\");\nstartRange.insertText(\"{m:enduserdoc}\", \"Start\");\nawait context.sync();\n\n// Step 2: remove the now-orphaned field (fldChar begin / instrText /\n// fldChar end runs) -- its displayed code has already been replaced by\n// the literal text inserted above.\nconst fields = context.document.body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nfields.items[0].delete();\nawait context.sync();\n\n// Step 3: prefix the visible error text with the \"    <---\" marker.\nconst matches = body.search(\"Invalid block: Unexpected tag m:enduserdoc at this location\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nmatches.items[0].insertText(\"    <---Invalid block: Unexpected tag m:enduserdoc at this location\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The paragraph currently contains a Word field \"{ m:enduserdoc }\" (a\n# fldChar begin / instrText / fldChar end run triple) followed by a bold\n# red run with an \"unexpected tag\" error message. The M2Doc parser was\n# switched to a token-iterator based field rewriter, so the expected\n# generation output now shows the raw field code as literal template\n# text (\"{m:enduserdoc}\") instead of an actual Word field, and the error\n# message gets a leading arrow marker (\"    <---\").\n\n# Step 1: insert the literal text \"{m:enduserdoc}\" right before the\n# field's code range. Inserting here (rather than after removing the\n# field) means the new run sits next to the field's begin run, so it\n# does not inherit the run properties (bold / red) of the trailing\n# error-message run.\n$f = $d.Fields.Item(1)\n$code = $f.Code\n$insertPoint = $d.Range($code.Start, $code.Start)\n$insertPoint.InsertBefore(\"{m:enduserdoc}\")\n\n# Step 2: remove the now-orphaned field (fldChar begin / instrText /\n# fldChar end runs) -- the field's displayed code has already been\n# replaced by the literal text inserted above.\n$f2 = $d.Fields.Item(1)\n$f2.Delete()\n\n# Step 3: prefix the visible error text with the \"    <---\" marker.\n$find = $d.Content.Find\n$find.Execute(\"Invalid block: Unexpected tag m:enduserdoc at this location\", $false, $false, $false, $false, $false, $true, 1, $false, \"    <---Invalid block: Unexpected tag m:enduserdoc at this location\", 2)\n"}
